$d = $word.ActiveDocument

# 1) Update the letter date: "September 19, 2025" -> "September 21, 2025"
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2) Split the mailing-address line "999 Story Road, San Jose CA 95122" (the one
#    that is NOT inside the property-address table) into two paragraphs:
#       "999 Story Road"
#       "San Jose, CA 95122"
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*999 Story Road, San Jose CA 95122*" -and $p.Range.Information(12) -eq $false) {
        $p.Range.Find.Execute("999 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false, $true, 1, $false, "999 Story Road", 2) | Out-Null
        break
    }
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*999 Story Road*" -and $p.Range.Information(12) -eq $false) {
        $p.Range.InsertParagraphAfter() | Out-Null
        $p2 = $p.Next()
        $p2.Range.Text = "San Jose, CA 95122"
        break
    }
}

# 3) Remove the now-superfluous empty "NoSpacing" paragraph that followed
#    "Board of Directors".
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Board of Directors*") {
        $nextP = $p.Next()
        $nextText = $nextP.Range.Text
        $nextStyle = $nextP.Style.NameLocal
        if ($nextText.Trim() -eq "" -and $nextStyle -eq "No Spacing") {
            $nextP.Range.Delete() | Out-Null
        }
        break
    }
}
